$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# Mark rows 262..295 (inclusive) as translated ("ok") in column B,
# matching the style already used throughout column B.
for ($r = 262; $r -le 295; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = "ok"
    $cell.HorizontalAlignment = -4108
}

# Restore the top selection, then move the split/freeze pane's
# top-left cell and the active selection down to follow the newly
# filled rows (topLeftCell A292, active cell B295 in the split pane).
$ws.Range("D1:E4").Select() | Out-Null
$excel.ActiveWindow.Split = $false
$ws.Range("A292").Select() | Out-Null
$excel.ActiveWindow.Split = $true
$ws.Range("B295").Select() | Out-Null
